$d = $word.ActiveDocument

# --- 1. Normalize the split "I" / "s the standardized..." runs into one run.
# (A no-op replace of the already-correct text forces Word to re-emit the
# sentence as a single run, same as the author's underlying edit.)
$d.Content.Find.Execute(
    "Is the standardized software coding style (for Python) being adhered?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Is the standardized software coding style (for Python) being adhered?", 2) | Out-Null

# --- 2. Reduce the Viva total from 10 to 8 points.
$d.Content.Find.Execute(
    "Viva [10 points]", $true, $false, $false, $false, $false, $true, 1, $false,
    "Viva [8 points]", 2) | Out-Null

# --- 3. Reduce the "Ability to answer questions from the Assignments 1-10." sub score from [2] to [1].
$pAssignments = $d.Paragraphs.Item(68)
$pAssignments.Range.Find.Execute(
    "[2]", $true, $false, $false, $false, $false, $true, 1, $false, "[1]", 2) | Out-Null

# --- 4. Reduce the "Ability to answer questions about the MLOps tools and usecases." sub score from [2] to [1].
$pMlops = $d.Paragraphs.Item(69)
$pMlops.Range.Find.Execute(
    " [2]", $true, $false, $false, $false, $false, $true, 1, $false, " [1]", 2) | Out-Null

# The author's cursor ended up right after the "1" (before the closing "]")
# when they typed the replacement digit; that's where the _GoBack bookmark
# now lives. Re-seat the bookmark there (Bookmarks.Add moves it if it
# already exists elsewhere, same as real Word).
$pMlopsEnd = $pMlops.Range.End
$bmRange = $d.Range($pMlopsEnd - 1, $pMlopsEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 5. The empty list paragraph right after the "duo project" bullet picks
# up an explicit ilvl=0 (outdenting it once re-homes it at level 0 while
# keeping numId=0).
$pEmpty = $d.Paragraphs.Item(76)
$pEmpty.Range.ListFormat.ListOutdent() | Out-Null

# --- 6. Flag "Default Paragraph Font" as a quick style (adds <w:qFormat/>).
$defaultParaFont = $d.Styles("Default Paragraph Font")
$defaultParaFont.QuickStyle = $true

Write-Output "done"
